$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing shared strings to reuse (read via Value2 to avoid COM ParameterizedProperty quirk with Value getter)
$catDef = $ws.Range("B72").Value2
$dogDef = $ws.Range("B73").Value2
$cityDef = $ws.Range("B6").Value2
$toyDef = $ws.Range("B10").Value2

$story = 'story'
$storyDef = @'
@story /'stɔ:ri/
*  danh từ
- chuyện, câu chuyện
=they all tell the same story+ họ đều kể một câu chuyện như nhau
=as the story goes+ người ta nói chuyện rằng
=but that is another story+ nhưng đó lại là chuyện khác
- truyện
=a short story+ truyện ngắn
- cốt truyện, tình tiết (một truyện, một vở kịch...)
=he reads only for the story+ anh ta đọc để hiểu cốt truyện thôi
- tiểu sử, quá khứ (của một người)
- luây kàng ngốc khoành người nói dối
=oh you story!+ nói dối!, điêu!
- (từ cổ,nghĩa cổ) lịch sử, sử học
*  danh từ
- (như) storey

'@
$fish = 'fish'
$fishDef = @'
@fish /fiʃ/
*  danh từ
- cá
=freshwater fish+ cá nước ngọt
=salt-water fish+ cá nước mặn
- cá, món cá
- (thiên văn học) chòm sao Cá
- người cắn câu, người bị mồi chài
- con người gã (có cá tính đặc biệt)
=a queer fish+ một con người (gã) kỳ quặc
!all's fish that comes to his net
- lớn bé, to nhỏ hắn quơ tất
!to be as drunk as a fish
- say bí tỉ
!to be as mute as a fish
- câm như hến
!to drink like a fish
- (xem) drink
!to feed the fishes
- chết đuối
- bị say sóng
!like a fish out of water
- (xem) water
!to have other fish to fly
- có công việc khác quan trọng hơn
!he who would catch fish must not mind getting wet
- muốn ăn cá phải lội nước, muốn ăn hét phải đào giun
!neither fish, fish, not good red herring
- môn chẳng ra môn, khoai chẳng ra khoai
!never fry a fish till it's caught
- chưa làm vòng chớ vội mong ăn thịt
!never offer to teach fish to swim
- chớ nên múa rìu qua mắt thợ
!a pretty kettle of fish
- (xem) kettle
!there's as good fish in the sea as ever came out of it
- thừa mứa chứa chan, nhiều vô kể
*  nội động từ
- đánh cá, câu cá, bắt cá
=to fish in the sea+ đánh cá ở biển
- (+ for) tìm, mò (cái gì ở dưới nước)
- (+ for) câu, moi những điều bí mật
*  ngoại động từ
- câu cá ở, đánh cá ở, bắt cá ở
=to fish a river+ đánh cá ở sông
- (hàng hải)
=to fish the anchor+ nhổ neo
- rút, lấy, kéo, moi
=to fish something out of water+ kéo cái gì từ dưới nước lên
- (từ hiếm,nghĩa hiếm) câu (cá), đánh (cá), bắt (cá), tìm (san hô...)
=to fish a troud+ câu một con cá hồi
!to fish out
- đánh hết cá (ở ao...)
- moi (ý kiến, bí mật)
!to fish in troubled waters
- lợi dụng đục nước béo cò
*  danh từ
- (hàng hải) miếng gỗ nẹp, miếng sắt nẹp (ở cột buồm, ở chỗ nối)
- (ngành đường sắt) thanh nối ray ((cũng) fish plate)
*  ngoại động từ
- (hàng hải) nẹp (bằng gỗ hay sắt)
- nối (đường ray) bằng thanh nối ray
*  danh từ
- (đánh bài) thẻ (bằng ngà... dùng thay tiền để đánh bài)

'@
$jsonCat1 = '{"code":300,"content":"cat","timestamp":"2020-07-09 18:34:00","dest":"0.0.0.1:8080","src":"0.0.0.1:64244"}'
$jsonDog1 = '{"code":300,"content":"dog","timestamp":"2020-07-09 18:34:09","dest":"0.0.0.1:8080","src":"0.0.0.1:64244"}'
$jsonCity1 = '{"code":300,"content":"city","timestamp":"2020-07-09 18:34:51","dest":"0.0.0.1:8080","src":"0.0.0.1:64244"}'
$jsonBig1 = '{"code":300,"content":"big","timestamp":"2020-07-09 18:35:02","dest":"0.0.0.1:8080","src":"0.0.0.1:64244"}'
$bigDef = @'
@big /big/
*  tính từ
- to, lớn
=a big tree+ cây to
=big repair+ sửa chữa lớn
=Big Three+ ba nước lớn
=Big Five+ năm nước lớn
- bụng to, có mang, có chửa
=big with news+ đầy tin, nhiều tin
- quan trọng
=a big man+ nhân vật quan trọng
- hào hiệp, phóng khoáng, rộng lượng
=he has a big hear+ anh ta là người hào hiệp
- huênh hoang, khoác lác
=big words+ những lời nói huênh hoang khoác lác
=big words+ những lời nói huênh hoang
!too big for one's boots (breeches, shoes, trousers)
- (từ lóng) quá tự tin, tự phụ tự mãn; làm bộ làm tịch
*  phó từ
- ra vẻ quan trọng, với vẻ quan trọng
=to look big+ làm ra vẻ quan trọng
- huênh hoang khoác lác
=to talk big+ nói huênh hoang, nói phách

'@
$jsonCity2 = '{"code":300,"content":"city","timestamp":"2020-07-09 18:35:06","dest":"0.0.0.1:8080","src":"0.0.0.1:64244"}'
$jsonCat2 = '{"code":300,"content":"cat","timestamp":"2020-07-09 18:37:47","dest":"0.0.0.1:8080","src":"0.0.0.1:64244"}'
$jsonToy1 = '{"code":300,"content":"toy","timestamp":"2020-07-09 18:38:01","dest":"0.0.0.1:8080","src":"0.0.0.1:64445"}'
$jsonCat3 = '{"code":300,"content":"cat","timestamp":"2020-07-10 01:13:29","dest":"0.0.0.1:8080","src":"0.0.0.1:60893"}'
$jsonSent1 = '{"code":300,"content":"sent","timestamp":"2020-07-10 01:13:42","dest":"0.0.0.1:8080","src":"0.0.0.1:60893"}'
$sentDef = @'
@sent /send/
*  ngoại động từ sent 
/sent/
- gửi, sai, phái, cho đi ((cũng) scend)
=to send word to somebody+ gửi vài chữ cho ai
=to send a boy a school+ cho một em nhỏ đi học
- cho, ban cho, phù hộ cho, giáng (trời,  Thượng đế...)
=send him victorioussend+ trời phụ hộ cho nó thắng trận!
=to send a drought+ giáng xuống nạn hạn hán
- bắn ra, làm bốc lên, làm nẩy ra, toả ra
=to send a ball over the trees+ đá tung quả bóng qua rặng cây
=to send smoke high in the air+ làm bốc khói lên cao trong không trung
- đuổi đi, tống đi
=to send somebody about his business+ tống cổ ai đi
- làm cho (mê mẩn)
=to send somebody crazy+ làm ai say mê; (nhạc ja, lóng) làm cho mê li
- (từ Mỹ,nghĩa Mỹ) hướng tới, đẩy tới
=your question has sent me to the dictionary+ câu hỏi của anh đã khiến tôi đi tìm từ điển
*  nội động từ
- gửi thư, nhắn
=to send to worn somebody+ gửi thư báo cho ai; gửi thư cảnh cáo ai
=to send to somebody to take care+ nhắn ai phải cẩn thận
!to send away
- gửi đi
- đuổi di
!to send after
- cho đi tìm, cho đuổi theo
!to send down
- cho xuống
- tạm đuổi, đuổi (khỏi trường)
!to send for
- gửi đặt mua
=to send for something+ gửi đặt mua cái gì
- nhắn đến, cho đi tìm đến, cho mời đến
=to send for somebody+ nhắn ai tìm đến, cho người mời đến
=to send for somebody+ nhắn ai đến, cho người mời ai
!to send forth
- toả ra, bốc ra (hương thơm, mùi, khói...)
- nảy ra (lộc non, lá...)
!to send in
- nộp, giao (đơn từ...)
- ghi, đăng (tên...)
=to send in one's name+ đăng tên (ở kỳ thi)
!to send off
- gửi đi (thư, quà) phái (ai) đi (công tác)
- đuổi đi, tống khứ
- tiễn đưa, hoan tống
!to send out
- gửi đi, phân phát
- toả ra, bốc ra (hương thơm, mùi, khói...)
- nảy ra
=trees send out young leaves+ cây ra lá non
!to send round
- chuyền tay, chuyền vòng (vật gì)
!to send up
- làm đứng dậy, làm trèo lên
- (từ Mỹ,nghĩa Mỹ),  (thông tục) kết án tù
!to send coals to Newcastle
- (xem) coal
!to send flying
- đuổi đi, bắt hối hả ra đi
- làm cho lảo đảo, đánh bật ra (bằng một cái đòn)
- làm cho chạy tan tác, làm cho tan tác
!to send packing
- đuổi đi, tống cổ đi
!to send someone to Jericho
- đuổi ai đi, tống cổ ai đi
!to send to Coventry
- phớt lờ, không hợp tác với (ai)

'@
$jsonSent2 = '{"code":300,"content":"sent","timestamp":"2020-07-10 01:13:47","dest":"0.0.0.1:8080","src":"0.0.0.1:60893"}'
$jsonCat4 = '{"code":300,"content":"cat","timestamp":"2020-07-10 15:31:39","dest":"0.0.0.1:8080","src":"0.0.0.1:58725"}'
$jsonToy2 = '{"code":300,"content":"toy","timestamp":"2020-07-10 15:31:55","dest":"0.0.0.1:8080","src":"0.0.0.1:58731"}'
$jsonFish1 = '{"code":300,"content":"fish","timestamp":"2020-07-10 15:32:27","dest":"0.0.0.1:8080","src":"0.0.0.1:58725"}'

$rows = @(
    @{ Row = 86;  A = $story;     B = $storyDef },
    @{ Row = 87;  A = $fish;      B = $fishDef },
    @{ Row = 88;  A = "cat";      B = $catDef },
    @{ Row = 89;  A = "dog";      B = $dogDef },
    @{ Row = 90;  A = "cat";      B = $catDef },
    @{ Row = 91;  A = "cat";      B = $catDef },
    @{ Row = 92;  A = $jsonCat1;  B = $catDef },
    @{ Row = 93;  A = $jsonDog1;  B = $dogDef },
    @{ Row = 94;  A = $jsonCity1; B = $cityDef },
    @{ Row = 95;  A = $jsonBig1;  B = $bigDef },
    @{ Row = 96;  A = $jsonCity2; B = $cityDef },
    @{ Row = 97;  A = $jsonCat2;  B = $catDef },
    @{ Row = 98;  A = $jsonToy1;  B = $toyDef },
    @{ Row = 99;  A = $jsonCat3;  B = $catDef },
    @{ Row = 100; A = $jsonSent1; B = $sentDef },
    @{ Row = 101; A = $jsonSent2; B = $sentDef },
    @{ Row = 102; A = $jsonCat4;  B = $catDef },
    @{ Row = 103; A = $jsonToy2;  B = $toyDef },
    @{ Row = 104; A = $jsonFish1; B = $fishDef }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.A
    $ws.Cells.Item($r.Row, 2).Value2 = $r.B
}
